$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 133
$ws.Range("D133").Value = $true
$ws.Range("E133").Value = $true
$ws.Range("H133").Value = 90
$ws.Range("J133").Value = $true
$ws.Range("K133").Value = '[{''identifier'': ''goals_scored'', ''a'': [{''value'': 1, ''element'': 214}, {''value'': 1, ''element'': 215}], ''h'': [{''value'': 1, ''element'': 269}, {''value'': 1, ''element'': 498}]}, {''identifier'': ''assists'', ''a'': [{''value'': 1, ''element'': 219}], ''h'': [{''value'': 1, ''element'': 266}, {''value'': 1, ''element'': 267}]}, {''identifier'': ''own_goals'', ''a'': [], ''h'': []}, {''identifier'': ''penalties_saved'', ''a'': [], ''h'': []}, {''identifier'': ''penalties_missed'', ''a'': [], ''h'': []}, {''identifier'': ''yellow_cards'', ''a'': [{''value'': 1, ''element'': 221}, {''value'': 1, ''element'': 222}], ''h'': [{''value'': 1, ''element'': 256}]}, {''identifier'': ''red_cards'', ''a'': [], ''h'': []}, {''identifier'': ''saves'', ''a'': [{''value'': 1, ''element'': 212}], ''h'': [{''value'': 7, ''element'': 262}]}, {''identifier'': ''bonus'', ''a'': [{''value'': 3, ''element'': 215}], ''h'': [{''value'': 2, ''element'': 269}, {''value'': 1, ''element'': 262}]}, {''identifier'': ''bps'', ''a'': [{''value'': 35, ''element'': 215}, {''value'': 21, ''element'': 219}, {''value'': 20, ''element'': 214}, {''value'': 17, ''element'': 204}, {''value'': 16, ''element'': 207}, {''value'': 15, ''element'': 203}, {''value'': 15, ''element'': 217}, {''value'': 13, ''element'': 221}, {''value'': 12, ''element'': 222}, {''value'': 11, ''element'': 212}, {''value'': 3, ''element'': 220}, {''value'': 3, ''element'': 443}, {''value'': 1, ''element'': 211}, {''value'': -2, ''element'': 218}], ''h'': [{''value'': 27, ''element'': 269}, {''value'': 24, ''element'': 262}, {''value'': 22, ''element'': 498}, {''value'': 17, ''element'': 254}, {''value'': 17, ''element'': 266}, {''value'': 15, ''element'': 259}, {''value'': 14, ''element'': 255}, {''value'': 13, ''element'': 267}, {''value'': 11, ''element'': 500}, {''value'': 8, ''element'': 271}, {''value'': 5, ''element'': 256}, {''value'': 5, ''element'': 466}, {''value'': 2, ''element'': 554}]}]'
$ws.Range("N133").Value = 2
$ws.Range("Q133").Value = 2

# Row 134
$ws.Range("D134").Value = $true
$ws.Range("E134").Value = $true
$ws.Range("H134").Value = 90
$ws.Range("J134").Value = $true
$ws.Range("K134").Value = '[{''identifier'': ''goals_scored'', ''a'': [{''value'': 1, ''element'': 133}, {''value'': 1, ''element'': 137}], ''h'': []}, {''identifier'': ''assists'', ''a'': [{''value'': 1, ''element'': 127}], ''h'': []}, {''identifier'': ''own_goals'', ''a'': [], ''h'': []}, {''identifier'': ''penalties_saved'', ''a'': [], ''h'': []}, {''identifier'': ''penalties_missed'', ''a'': [], ''h'': []}, {''identifier'': ''yellow_cards'', ''a'': [], ''h'': [{''value'': 1, ''element'': 86}, {''value'': 1, ''element'': 97}]}, {''identifier'': ''red_cards'', ''a'': [], ''h'': []}, {''identifier'': ''saves'', ''a'': [{''value'': 3, ''element'': 131}], ''h'': [{''value'': 1, ''element'': 93}]}, {''identifier'': ''bonus'', ''a'': [{''value'': 3, ''element'': 133}, {''value'': 2, ''element'': 124}, {''value'': 1, ''element'': 131}], ''h'': []}, {''identifier'': ''bps'', ''a'': [{''value'': 33, ''element'': 133}, {''value'': 29, ''element'': 124}, {''value'': 28, ''element'': 131}, {''value'': 24, ''element'': 128}, {''value'': 23, ''element'': 123}, {''value'': 22, ''element'': 127}, {''value'': 20, ''element'': 137}, {''value'': 12, ''element'': 135}, {''value'': 10, ''element'': 138}, {''value'': 9, ''element'': 134}, {''value'': 5, ''element'': 125}, {''value'': 5, ''element'': 468}, {''value'': 1, ''element'': 129}, {''value'': -2, ''element'': 139}], ''h'': [{''value'': 18, ''element'': 100}, {''value'': 16, ''element'': 98}, {''value'': 13, ''element'': 84}, {''value'': 11, ''element'': 83}, {''value'': 11, ''element'': 93}, {''value'': 9, ''element'': 86}, {''value'': 9, ''element'': 102}, {''value'': 8, ''element'': 97}, {''value'': 5, ''element'': 430}, {''value'': 3, ''element'': 87}, {''value'': 3, ''element'': 451}, {''value'': 1, ''element'': 90}, {''value'': -1, ''element'': 447}, {''value'': -5, ''element'': 91}]}]'
$ws.Range("N134").Value = 2
$ws.Range("Q134").Value = 0

# Row 135
$ws.Range("D135").Value = $true
$ws.Range("E135").Value = $true
$ws.Range("H135").Value = 90
$ws.Range("J135").Value = $true
$ws.Range("K135").Value = '[{''identifier'': ''goals_scored'', ''a'': [{''value'': 1, ''element'': 376}], ''h'': []}, {''identifier'': ''assists'', ''a'': [{''value'': 1, ''element'': 399}], ''h'': []}, {''identifier'': ''own_goals'', ''a'': [], ''h'': []}, {''identifier'': ''penalties_saved'', ''a'': [], ''h'': []}, {''identifier'': ''penalties_missed'', ''a'': [], ''h'': []}, {''identifier'': ''yellow_cards'', ''a'': [{''value'': 1, ''element'': 376}, {''value'': 1, ''element'': 392}, {''value'': 1, ''element'': 399}], ''h'': []}, {''identifier'': ''red_cards'', ''a'': [], ''h'': []}, {''identifier'': ''saves'', ''a'': [{''value'': 6, ''element'': 553}], ''h'': [{''value'': 3, ''element'': 111}]}, {''identifier'': ''bonus'', ''a'': [{''value'': 3, ''element'': 376}, {''value'': 2, ''element'': 553}, {''value'': 1, ''element'': 382}], ''h'': []}, {''identifier'': ''bps'', ''a'': [{''value'': 36, ''element'': 376}, {''value'': 34, ''element'': 553}, {''value'': 30, ''element'': 382}, {''value'': 26, ''element'': 377}, {''value'': 24, ''element'': 380}, {''value'': 15, ''element'': 393}, {''value'': 15, ''element'': 399}, {''value'': 14, ''element'': 395}, {''value'': 9, ''element'': 390}, {''value'': 8, ''element'': 396}, {''value'': 3, ''element'': 457}, {''value'': 2, ''element'': 378}, {''value'': -1, ''element'': 392}, {''value'': -2, ''element'': 389}], ''h'': [{''value'': 22, ''element'': 506}, {''value'': 19, ''element'': 107}, {''value'': 17, ''element'': 432}, {''value'': 15, ''element'': 459}, {''value'': 14, ''element'': 463}, {''value'': 13, ''element'': 111}, {''value'': 10, ''element'': 114}, {''value'': 9, ''element'': 118}, {''value'': 9, ''element'': 529}, {''value'': 7, ''element'': 119}, {''value'': 6, ''element'': 113}, {''value'': 6, ''element'': 116}, {''value'': 1, ''element'': 110}, {''value'': -1, ''element'': 431}]}]'
$ws.Range("N135").Value = 1
$ws.Range("Q135").Value = 0

# Row 136
$ws.Range("D136").Value = $true
$ws.Range("E136").Value = $true
$ws.Range("H136").Value = 90
$ws.Range("J136").Value = $true
$ws.Range("K136").Value = '[{''identifier'': ''goals_scored'', ''a'': [{''value'': 1, ''element'': 42}], ''h'': [{''value'': 2, ''element'': 183}]}, {''identifier'': ''assists'', ''a'': [{''value'': 1, ''element'': 421}], ''h'': [{''value'': 2, ''element'': 182}]}, {''identifier'': ''own_goals'', ''a'': [], ''h'': []}, {''identifier'': ''penalties_saved'', ''a'': [], ''h'': []}, {''identifier'': ''penalties_missed'', ''a'': [], ''h'': []}, {''identifier'': ''yellow_cards'', ''a'': [], ''h'': []}, {''identifier'': ''red_cards'', ''a'': [], ''h'': [{''value'': 1, ''element'': 189}]}, {''identifier'': ''saves'', ''a'': [{''value'': 5, ''element'': 47}], ''h'': [{''value'': 2, ''element'': 189}, {''value'': 2, ''element'': 526}]}, {''identifier'': ''bonus'', ''a'': [{''value'': 1, ''element'': 42}], ''h'': [{''value'': 3, ''element'': 183}, {''value'': 2, ''element'': 182}]}, {''identifier'': ''bps'', ''a'': [{''value'': 27, ''element'': 42}, {''value'': 24, ''element'': 475}, {''value'': 19, ''element'': 47}, {''value'': 17, ''element'': 54}, {''value'': 17, ''element'': 57}, {''value'': 16, ''element'': 493}, {''value'': 16, ''element'': 516}, {''value'': 11, ''element'': 49}, {''value'': 9, ''element'': 38}, {''value'': 9, ''element'': 53}, {''value'': 7, ''element'': 421}, {''value'': 6, ''element'': 537}, {''value'': 1, ''element'': 502}, {''value'': -1, ''element'': 534}], ''h'': [{''value'': 42, ''element'': 183}, {''value'': 29, ''element'': 182}, {''value'': 17, ''element'': 192}, {''value'': 15, ''element'': 186}, {''value'': 13, ''element'': 199}, {''value'': 12, ''element'': 198}, {''value'': 9, ''element'': 181}, {''value'': 7, ''element'': 189}, {''value'': 7, ''element'': 193}, {''value'': 6, ''element'': 191}, {''value'': 5, ''element'': 187}, {''value'': 5, ''element'': 526}, {''value'': 3, ''element'': 188}, {''value'': 3, ''element'': 195}]}]'
$ws.Range("N136").Value = 1
$ws.Range("Q136").Value = 2

# Row 137
$ws.Range("D137").Value = $true
$ws.Range("E137").Value = $true
$ws.Range("H137").Value = 90
$ws.Range("J137").Value = $true
$ws.Range("K137").Value = '[{''identifier'': ''goals_scored'', ''a'': [{''value'': 2, ''element'': 505}], ''h'': [{''value'': 2, ''element'': 344}, {''value'': 1, ''element'': 349}]}, {''identifier'': ''assists'', ''a'': [{''value'': 1, ''element'': 79}, {''value'': 1, ''element'': 492}], ''h'': [{''value'': 2, ''element'': 342}, {''value'': 1, ''element'': 331}]}, {''identifier'': ''own_goals'', ''a'': [], ''h'': []}, {''identifier'': ''penalties_saved'', ''a'': [], ''h'': []}, {''identifier'': ''penalties_missed'', ''a'': [], ''h'': []}, {''identifier'': ''yellow_cards'', ''a'': [{''value'': 1, ''element'': 65}, {''value'': 1, ''element'': 82}], ''h'': []}, {''identifier'': ''red_cards'', ''a'': [], ''h'': []}, {''identifier'': ''saves'', ''a'': [{''value'': 2, ''element'': 494}], ''h'': [{''value'': 4, ''element'': 341}]}, {''identifier'': ''bonus'', ''a'': [{''value'': 1, ''element'': 505}], ''h'': [{''value'': 3, ''element'': 344}, {''value'': 2, ''element'': 342}]}, {''identifier'': ''bps'', ''a'': [{''value'': 38, ''element'': 505}, {''value'': 29, ''element'': 79}, {''value'': 24, ''element'': 492}, {''value'': 14, ''element'': 58}, {''value'': 14, ''element'': 59}, {''value'': 11, ''element'': 494}, {''value'': 10, ''element'': 82}, {''value'': 10, ''element'': 446}, {''value'': 9, ''element'': 69}, {''value'': 7, ''element'': 65}, {''value'': 5, ''element'': 75}, {''value'': 1, ''element'': 67}], ''h'': [{''value'': 50, ''element'': 344}, {''value'': 47, ''element'': 342}, {''value'': 36, ''element'': 349}, {''value'': 26, ''element'': 331}, {''value'': 19, ''element'': 442}, {''value'': 17, ''element'': 341}, {''value'': 16, ''element'': 336}, {''value'': 15, ''element'': 330}, {''value'': 14, ''element'': 333}, {''value'': 14, ''element'': 350}, {''value'': 6, ''element'': 338}, {''value'': 4, ''element'': 345}, {''value'': 3, ''element'': 347}, {''value'': 2, ''element'': 523}]}]'
$ws.Range("N137").Value = 2
$ws.Range("Q137").Value = 3

# Row 138
$ws.Range("D138").Value = $true
$ws.Range("E138").Value = $true
$ws.Range("H138").Value = 90
$ws.Range("J138").Value = $true
$ws.Range("K138").Value = '[{''identifier'': ''goals_scored'', ''a'': [{''value'': 1, ''element'': 525}], ''h'': [{''value'': 1, ''element'': 313}, {''value'': 1, ''element'': 321}]}, {''identifier'': ''assists'', ''a'': [{''value'': 1, ''element'': 373}], ''h'': [{''value'': 1, ''element'': 325}, {''value'': 1, ''element'': 438}]}, {''identifier'': ''own_goals'', ''a'': [], ''h'': []}, {''identifier'': ''penalties_saved'', ''a'': [], ''h'': []}, {''identifier'': ''penalties_missed'', ''a'': [], ''h'': []}, {''identifier'': ''yellow_cards'', ''a'': [], ''h'': [{''value'': 1, ''element'': 320}]}, {''identifier'': ''red_cards'', ''a'': [], ''h'': []}, {''identifier'': ''saves'', ''a'': [{''value'': 3, ''element'': 366}], ''h'': [{''value'': 2, ''element'': 318}]}, {''identifier'': ''bonus'', ''a'': [{''value'': 1, ''element'': 373}], ''h'': [{''value'': 3, ''element'': 313}, {''value'': 2, ''element'': 321}]}, {''identifier'': ''bps'', ''a'': [{''value'': 28, ''element'': 373}, {''value'': 21, ''element'': 525}, {''value'': 16, ''element'': 357}, {''value'': 15, ''element'': 366}, {''value'': 12, ''element'': 358}, {''value'': 11, ''element'': 356}, {''value'': 10, ''element'': 370}, {''value'': 9, ''element'': 355}, {''value'': 8, ''element'': 353}, {''value'': 8, ''element'': 368}, {''value'': 4, ''element'': 542}, {''value'': 3, ''element'': 363}, {''value'': 3, ''element'': 374}, {''value'': -1, ''element'': 364}], ''h'': [{''value'': 38, ''element'': 313}, {''value'': 32, ''element'': 321}, {''value'': 27, ''element'': 325}, {''value'': 24, ''element'': 438}, {''value'': 17, ''element'': 486}, {''value'': 14, ''element'': 307}, {''value'': 13, ''element'': 318}, {''value'': 10, ''element'': 308}, {''value'': 9, ''element'': 305}, {''value'': 8, ''element'': 508}, {''value'': 6, ''element'': 306}, {''value'': 4, ''element'': 314}, {''value'': -1, ''element'': 320}, {''value'': -2, ''element'': 315}]}]'
$ws.Range("N138").Value = 1
$ws.Range("Q138").Value = 2

# Row 139
$ws.Range("D139").Value = $true
$ws.Range("E139").Value = $true
$ws.Range("H139").Value = 90
$ws.Range("J139").Value = $true
$ws.Range("K139").Value = '[{''identifier'': ''goals_scored'', ''a'': [{''value'': 2, ''element'': 11}], ''h'': [{''value'': 1, ''element'': 278}, {''value'': 1, ''element'': 287}]}, {''identifier'': ''assists'', ''a'': [{''value'': 1, ''element'': 1}, {''value'': 1, ''element'': 15}], ''h'': [{''value'': 1, ''element'': 284}, {''value'': 1, ''element'': 286}]}, {''identifier'': ''own_goals'', ''a'': [], ''h'': []}, {''identifier'': ''penalties_saved'', ''a'': [], ''h'': []}, {''identifier'': ''penalties_missed'', ''a'': [], ''h'': []}, {''identifier'': ''yellow_cards'', ''a'': [{''value'': 1, ''element'': 467}], ''h'': [{''value'': 1, ''element'': 280}, {''value'': 1, ''element'': 286}]}, {''identifier'': ''red_cards'', ''a'': [], ''h'': []}, {''identifier'': ''saves'', ''a'': [{''value'': 6, ''element'': 14}], ''h'': [{''value'': 4, ''element'': 280}]}, {''identifier'': ''bonus'', ''a'': [{''value'': 3, ''element'': 11}, {''value'': 1, ''element'': 15}], ''h'': [{''value'': 2, ''element'': 278}, {''value'': 1, ''element'': 287}]}, {''identifier'': ''bps'', ''a'': [{''value'': 58, ''element'': 11}, {''value'': 26, ''element'': 15}, {''value'': 24, ''element'': 1}, {''value'': 20, ''element'': 14}, {''value'': 16, ''element'': 18}, {''value'': 15, ''element'': 106}, {''value'': 14, ''element'': 3}, {''value'': 13, ''element'': 21}, {''value'': 9, ''element'': 490}, {''value'': 6, ''element'': 467}, {''value'': 4, ''element'': 19}, {''value'': 3, ''element'': 504}, {''value'': 3, ''element'': 541}, {''value'': 1, ''element'': 12}], ''h'': [{''value'': 31, ''element'': 278}, {''value'': 26, ''element'': 287}, {''value'': 19, ''element'': 284}, {''value'': 18, ''element'': 507}, {''value'': 15, ''element'': 273}, {''value'': 15, ''element'': 290}, {''value'': 13, ''element'': 275}, {''value'': 13, ''element'': 280}, {''value'': 12, ''element'': 286}, {''value'': 12, ''element'': 454}, {''value'': 11, ''element'': 274}, {''value'': 8, ''element'': 283}]}]'
$ws.Range("N139").Value = 2
$ws.Range("Q139").Value = 2

# Row 140
$ws.Range("D140").Value = $true
$ws.Range("E140").Value = $true
$ws.Range("H140").Value = 90
$ws.Range("J140").Value = $true
$ws.Range("K140").Value = '[{''identifier'': ''goals_scored'', ''a'': [{''value'': 1, ''element'': 70}], ''h'': [{''value'': 1, ''element'': 401}]}, {''identifier'': ''assists'', ''a'': [{''value'': 1, ''element'': 294}], ''h'': [{''value'': 1, ''element'': 409}]}, {''identifier'': ''own_goals'', ''a'': [], ''h'': []}, {''identifier'': ''penalties_saved'', ''a'': [], ''h'': []}, {''identifier'': ''penalties_missed'', ''a'': [], ''h'': []}, {''identifier'': ''yellow_cards'', ''a'': [{''value'': 1, ''element'': 291}, {''value'': 1, ''element'': 293}, {''value'': 1, ''element'': 294}, {''value'': 1, ''element'': 295}, {''value'': 1, ''element'': 301}], ''h'': [{''value'': 1, ''element'': 420}]}, {''identifier'': ''red_cards'', ''a'': [], ''h'': []}, {''identifier'': ''saves'', ''a'': [{''value'': 3, ''element'': 471}], ''h'': [{''value'': 3, ''element'': 411}]}, {''identifier'': ''bonus'', ''a'': [{''value'': 3, ''element'': 70}, {''value'': 1, ''element'': 302}], ''h'': [{''value'': 2, ''element'': 401}]}, {''identifier'': ''bps'', ''a'': [{''value'': 30, ''element'': 70}, {''value'': 24, ''element'': 302}, {''value'': 18, ''element'': 423}, {''value'': 16, ''element'': 471}, {''value'': 15, ''element'': 294}, {''value'': 14, ''element'': 301}, {''value'': 13, ''element'': 291}, {''value'': 13, ''element'': 297}, {''value'': 7, ''element'': 293}, {''value'': 6, ''element'': 501}, {''value'': 4, ''element'': 295}, {''value'': 4, ''element'': 453}, {''value'': 2, ''element'': 441}, {''value'': -1, ''element'': 303}], ''h'': [{''value'': 25, ''element'': 401}, {''value'': 18, ''element'': 415}, {''value'': 16, ''element'': 402}, {''value'': 16, ''element'': 403}, {''value'': 15, ''element'': 409}, {''value'': 14, ''element'': 411}, {''value'': 13, ''element'': 414}, {''value'': 9, ''element'': 408}, {''value'': 9, ''element'': 417}, {''value'': 8, ''element'': 420}, {''value'': 5, ''element'': 410}]}]'
$ws.Range("N140").Value = 1
$ws.Range("Q140").Value = 1

# Row 141
$ws.Range("D141").Value = $true
$ws.Range("E141").Value = $true
$ws.Range("H141").Value = 90
$ws.Range("J141").Value = $true
$ws.Range("K141").Value = '[{''identifier'': ''goals_scored'', ''a'': [{''value'': 1, ''element'': 150}], ''h'': [{''value'': 1, ''element'': 166}, {''value'': 1, ''element'': 167}]}, {''identifier'': ''assists'', ''a'': [{''value'': 1, ''element'': 511}], ''h'': [{''value'': 1, ''element'': 159}, {''value'': 1, ''element'': 167}]}, {''identifier'': ''own_goals'', ''a'': [], ''h'': []}, {''identifier'': ''penalties_saved'', ''a'': [], ''h'': []}, {''identifier'': ''penalties_missed'', ''a'': [], ''h'': []}, {''identifier'': ''yellow_cards'', ''a'': [{''value'': 1, ''element'': 151}], ''h'': []}, {''identifier'': ''red_cards'', ''a'': [], ''h'': []}, {''identifier'': ''saves'', ''a'': [{''value'': 4, ''element'': 148}], ''h'': [{''value'': 2, ''element'': 168}]}, {''identifier'': ''bonus'', ''a'': [{''value'': 1, ''element'': 150}], ''h'': [{''value'': 3, ''element'': 167}, {''value'': 2, ''element'': 159}]}, {''identifier'': ''bps'', ''a'': [{''value'': 26, ''element'': 150}, {''value'': 20, ''element'': 148}, {''value'': 19, ''element'': 511}, {''value'': 15, ''element'': 141}, {''value'': 14, ''element'': 142}, {''value'': 13, ''element'': 154}, {''value'': 12, ''element'': 17}, {''value'': 12, ''element'': 478}, {''value'': 8, ''element'': 151}, {''value'': 7, ''element'': 143}, {''value'': 4, ''element'': 147}, {''value'': 4, ''element'': 158}, {''value'': 1, ''element'': 496}], ''h'': [{''value'': 38, ''element'': 167}, {''value'': 27, ''element'': 159}, {''value'': 25, ''element'': 166}, {''value'': 21, ''element'': 162}, {''value'': 20, ''element'': 176}, {''value'': 16, ''element'': 164}, {''value'': 15, ''element'': 161}, {''value'': 13, ''element'': 168}, {''value'': 13, ''element'': 171}, {''value'': 13, ''element'': 448}, {''value'': 9, ''element'': 172}, {''value'': 3, ''element'': 174}, {''value'': 2, ''element'': 265}]}]'
$ws.Range("N141").Value = 1
$ws.Range("Q141").Value = 2

# Row 142
$ws.Range("D142").Value = $true
$ws.Range("E142").Value = $true
$ws.Range("H142").Value = 90
$ws.Range("J142").Value = $true
$ws.Range("K142").Value = '[{''identifier'': ''goals_scored'', ''a'': [{''value'': 1, ''element'': 29}, {''value'': 1, ''element'': 445}], ''h'': [{''value'': 1, ''element'': 225}]}, {''identifier'': ''assists'', ''a'': [{''value'': 1, ''element'': 30}, {''value'': 1, ''element'': 309}], ''h'': [{''value'': 1, ''element'': 233}]}, {''identifier'': ''own_goals'', ''a'': [{''value'': 1, ''element'': 94}], ''h'': []}, {''identifier'': ''penalties_saved'', ''a'': [], ''h'': []}, {''identifier'': ''penalties_missed'', ''a'': [], ''h'': []}, {''identifier'': ''yellow_cards'', ''a'': [{''value'': 1, ''element'': 458}], ''h'': [{''value'': 1, ''element'': 160}, {''value'': 1, ''element'': 226}, {''value'': 1, ''element'': 244}, {''value'': 1, ''element'': 549}]}, {''identifier'': ''red_cards'', ''a'': [], ''h'': []}, {''identifier'': ''saves'', ''a'': [{''value'': 5, ''element'': 94}], ''h'': [{''value'': 1, ''element'': 235}]}, {''identifier'': ''bonus'', ''a'': [{''value'': 3, ''element'': 309}, {''value'': 2, ''element'': 29}], ''h'': [{''value'': 1, ''element'': 225}]}, {''identifier'': ''bps'', ''a'': [{''value'': 31, ''element'': 309}, {''value'': 26, ''element'': 29}, {''value'': 23, ''element'': 445}, {''value'': 20, ''element'': 470}, {''value'': 17, ''element'': 458}, {''value'': 15, ''element'': 30}, {''value'': 15, ''element'': 94}, {''value'': 11, ''element'': 31}, {''value'': 8, ''element'': 28}, {''value'': 8, ''element'': 452}, {''value'': 6, ''element'': 465}, {''value'': 5, ''element'': 428}, {''value'': 1, ''element'': 35}], ''h'': [{''value'': 25, ''element'': 225}, {''value'': 14, ''element'': 246}, {''value'': 11, ''element'': 122}, {''value'': 11, ''element'': 549}, {''value'': 10, ''element'': 235}, {''value'': 10, ''element'': 244}, {''value'': 8, ''element'': 160}, {''value'': 7, ''element'': 233}, {''value'': 7, ''element'': 243}, {''value'': 4, ''element'': 241}, {''value'': 4, ''element'': 242}, {''value'': 3, ''element'': 234}, {''value'': 1, ''element'': 239}, {''value'': -3, ''element'': 226}]}]'
$ws.Range("N142").Value = 2
$ws.Range("Q142").Value = 2
